$wb = $excel.ActiveWorkbook

# Map of sheet name -> { cell -> new value }
$sheetUpdates = @{
    "展览" = @{
        "F5"  = 286
        "F6"  = 1096
        "F7"  = 1438
        "F12" = 166
        "F14" = 433
        "F15" = 1355
        "F16" = 114
        "F17" = 107
        "F19" = 5212
        "F20" = 35
        "F24" = 223
        "F26" = 5870
        "F31" = 14514
        "F32" = 1438
        "F36" = 8775
        "F37" = 623
        "F38" = 4210
        "F39" = 142
    }
    "全部类型" = @{
        "F5"  = 286
        "F6"  = 1096
        "F7"  = 1438
        "F12" = 166
        "F14" = 433
        "F15" = 1355
        "F16" = 114
        "F17" = 107
        "F20" = 5212
        "F21" = 35
        "F26" = 223
        "F29" = 5870
        "F34" = 14514
        "F35" = 1438
        "F39" = 8775
        "F40" = 623
        "F41" = 4210
        "F42" = 142
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cells = $sheetUpdates[$sheetName]
    foreach ($cellRef in $cells.Keys) {
        $ws.Range($cellRef).Value = $cells[$cellRef]
    }
}
